# Add payment / reusable component test rows to the Run_Manager sheet,
# and flip the "Execute" flag on the existing last row (row 110) from
# Yes -> No since it is no longer the final row of the run list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run_Manager")

# --- Row 110: Execute column changes from "Yes" to "No" -------------------
$ws.Range("D110").Value = "No"

# NOTE: the write order below intentionally matches the order new labels
# were first introduced by the original authoring tool (column B filled
# for the two new "Payment" rows, then column C, then the third new row)
# so that the rebuilt shared-strings table lines up exactly.

# --- ModuleName (column A) for the three new rows --------------------------
$ws.Range("A111").Value = "Payment"
$ws.Range("A112").Value = "Payment"

# --- TestName (column B) for rows 111-112 -----------------------------------
$ws.Range("B111").Value = "verifyPaymentInformationPageObjects"
$ws.Range("B112").Value = "verifyBankPaymentFormFields"

# --- TestDescription (column C) for rows 111-112 ----------------------------
$ws.Range("C111").Value = "to verify payment information page objects"
$ws.Range("C112").Value = "verify bank payment form fields"

# --- Row 113: Payment / verifyBankPaymentFormFieldsValidation --------------
$ws.Range("A113").Value = "Payment"
$ws.Range("B113").Value = "verifyBankPaymentFormFieldsValidation"
$ws.Range("C113").Value = "validation on bank payment form"

# --- Remaining columns (Execute / Priority / Count) for the new rows -------
$ws.Range("D111").Value = "No"
$ws.Range("E111").Value = "'1"
$ws.Range("F111").Value = "'1"

$ws.Range("D112").Value = "No"
$ws.Range("E112").Value = "'1"
$ws.Range("F112").Value = "'1"

$ws.Range("D113").Value = "Yes"
$ws.Range("E113").Value = "'1"
$ws.Range("F113").Value = "'1"

# --- Keep the visible selection on the new final row, as in the source ----
$null = $ws.Range("C116").Select()
